$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new formula/question pair for the ship/mass/temperature problem
$ws.Range("E3").Value = "(e*m)/Temp"
$ws.Range("F3").Value = "Das `$m(2,23) kg schwere schiff lädt `$e(3,7) Kisten bei `$Temp(25,30) °C auf."

# Update the view: scroll/zoom and selection moved from D3 to F3
$ws.Application.ActiveWindow.Zoom = 180
$ws.Range("F3").Select()
